# Apply updated symbol list values (cryptos.xlsx) - generated from diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''306.00'
$ws.Range("E2").Value = '''0.21%'
$ws.Range("D3").Value = '''40.33'
$ws.Range("E3").Value = '''2.07%'
$ws.Range("D4").Value = '''5.116'
$ws.Range("E4").Value = '''1.13%'
$ws.Range("E5").Value = '''-1.91%'
$ws.Range("E6").Value = '''-1.85%'
$ws.Range("D7").Value = '''0.9074'
$ws.Range("E7").Value = '''-0.86%'
$ws.Range("E8").Value = '''-6.08%'
$ws.Range("D9").Value = '''0.1011'
$ws.Range("E9").Value = '''4.30%'
$ws.Range("D10").Value = '''0.1754'
$ws.Range("E10").Value = '''1.42%'
$ws.Range("D11").Value = '''0.09157'
$ws.Range("E11").Value = '''2.86%'
$ws.Range("D12").Value = '''0.04177'
$ws.Range("E12").Value = '''-5.09%'
$ws.Range("D13").Value = '''0.1056'
$ws.Range("E13").Value = '''-0.23%'
$ws.Range("D14").Value = '''0.001246'
$ws.Range("E14").Value = '''-0.49%'
$ws.Range("D15").Value = '''0.005883'
$ws.Range("E15").Value = '''3.88%'
$ws.Range("D16").Value = '''3.353'
$ws.Range("E16").Value = '''-0.44%'
$ws.Range("D17").Value = '''4.267'
$ws.Range("E17").Value = '''-1.27%'
$ws.Range("D18").Value = '''0.3273'
$ws.Range("E18").Value = '''-2.77%'
$ws.Range("D19").Value = '''6.648'
$ws.Range("E19").Value = '''-5.69%'
$ws.Range("E20").Value = '''-0.69%'
$ws.Range("D21").Value = '''0.2727'
$ws.Range("E21").Value = '''-5.93%'
$ws.Range("D22").Value = '''0.04179'
$ws.Range("E22").Value = '''0.50%'
$ws.Range("E23").Value = '''1.56%'
$ws.Range("D24").Value = '''0.004060'
$ws.Range("E24").Value = '''-0.53%'
$ws.Range("D25").Value = '''0.0001303'
$ws.Range("E25").Value = '''5.70%'
$ws.Range("D26").Value = '''0.0003010'
$ws.Range("E26").Value = '''0.42%'
$ws.Range("E38").Value = '''1.19%'
$ws.Range("D39").Value = '''0.05164'
$ws.Range("E39").Value = '''0.54%'
$ws.Range("D40").Value = '''0.007780'
$ws.Range("E40").Value = '''-2.71%'
$ws.Range("D41").Value = '''0.1296'
$ws.Range("E41").Value = '''-2.38%'
$ws.Range("D42").Value = '''0.007068'
$ws.Range("E42").Value = '''-5.08%'
$ws.Range("D43").Value = '''0.001922'
$ws.Range("E43").Value = '''-6.18%'
$ws.Range("D44").Value = '''0.008440'
$ws.Range("E44").Value = '''5.11%'
$ws.Range("D45").Value = '''0.3305'
$ws.Range("E45").Value = '''-0.14%'
$ws.Range("D46").Value = '''0.00006361'
$ws.Range("E46").Value = '''-5.52%'
$ws.Range("E47").Value = '''-0.52%'
$ws.Range("B48").Value = 'CoinbaseStockToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D48").Value = '''0.004404'
$ws.Range("E48").Value = '''6.76%'
$ws.Range("B49").Value = 'BOLO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D49").Value = '''0.007105'
$ws.Range("E49").Value = '''106.65%'
$ws.Range("D50").Value = '''0.00002102'
$ws.Range("E50").Value = '''-0.52%'
$ws.Range("D51").Value = '''0.0002002'
$ws.Range("E51").Value = '''-0.52%'
